$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Order_Verongiida"
$ws.Cells.Item(2, 2).Value = -6.24020422541597
$ws.Cells.Item(2, 3).Value = [double]"0.0001024706194922948"
$ws.Cells.Item(2, 4).Value = "Impact`nlow_vs_medium"

# Row 3
$ws.Cells.Item(3, 1).Value = "Order_Verongiida"
$ws.Cells.Item(3, 2).Value = 8.839895652237658
$ws.Cells.Item(3, 3).Value = [double]"2.413483804509196E-06"
$ws.Cells.Item(3, 4).Value = "Impact`nmedium_vs_high"

# Row 4
$ws.Cells.Item(4, 1).Value = "Petrosiidae"
$ws.Cells.Item(4, 2).Value = -5.487336006413498
$ws.Cells.Item(4, 3).Value = [double]"1.782639053368633E-05"
$ws.Cells.Item(4, 4).Value = "Impact`nlow_vs_medium"

# Row 5
$ws.Cells.Item(5, 1).Value = "Chattonellaceae"
$ws.Cells.Item(5, 2).Value = -4.093055456223142
$ws.Cells.Item(5, 3).Value = [double]"0.0008718637146058619"
$ws.Cells.Item(5, 4).Value = "Impact`nlow_vs_medium"

# Row 6
$ws.Cells.Item(6, 1).Value = "Class_Phaeophyceae"
$ws.Cells.Item(6, 2).Value = -4.968706044194255
$ws.Cells.Item(6, 3).Value = [double]"4.254648742136053E-05"
$ws.Cells.Item(6, 4).Value = "Impact`nlow_vs_medium"

# Row 7
$ws.Cells.Item(7, 2).Value = 4.618869259448249
$ws.Cells.Item(7, 3).Value = [double]"0.0001199447140529064"

# Row 8
$ws.Cells.Item(8, 2).Value = 6.289912686711788
$ws.Cells.Item(8, 3).Value = [double]"1.845403527087003E-07"

# Row 9
$ws.Cells.Item(9, 2).Value = 4.112657550358105
$ws.Cells.Item(9, 3).Value = [double]"0.0007422452736906149"

# Row 10
$ws.Cells.Item(10, 2).Value = 6.549253003326004
$ws.Cells.Item(10, 3).Value = [double]"5.784284659792789E-07"

# Row 11
$ws.Cells.Item(11, 2).Value = 7.110415739267379
$ws.Cells.Item(11, 3).Value = [double]"5.11252100811555E-08"

# Row 12
$ws.Cells.Item(12, 2).Value = -4.951024750666711
$ws.Cells.Item(12, 3).Value = [double]"5.630960019763035E-05"

# Row 13
$ws.Cells.Item(13, 2).Value = 14.38946233558463
$ws.Cells.Item(13, 3).Value = [double]"0.0001385368952025626"

# Row 14
$ws.Cells.Item(14, 2).Value = 14.97160525940761
$ws.Cells.Item(14, 3).Value = [double]"0.0001381242253133509"
